# Revised ESST1A input file.
#
# Fills in the previously-blank ESST1A parameter header row (F1:Y1) on the
# "ESST1A" sheet with the model's field names, and moves the active
# selection on that sheet to R17 (with the view scrolled so column M is
# the first visible column after the frozen header row/column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESST1A")
[void]$ws.Activate()

$headers = [ordered]@{
    "F1" = "TR"
    "G1" = "VIMAX"
    "H1" = "VIMIN"
    "I1" = "TB"
    "J1" = "TC"
    "K1" = "TB1"
    "L1" = "TC1"
    "M1" = "VAMAX"
    "N1" = "VAMIN"
    "O1" = "KA"
    "P1" = "TA"
    "Q1" = "ILR"
    "R1" = "KLR"
    "S1" = "VRMAX"
    "T1" = "VRMIN"
    "U1" = "KF"
    "V1" = "TF"
    "W1" = "KC"
    "X1" = "UELc"
    "Y1" = "VOSc"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Scroll so column M is the first visible (unfrozen) column, and move the
# selection to R17 to match the saved view state.
try {
    $excel.ActiveWindow.ScrollColumn = 13
} catch {
}

[void]$ws.Range("R17").Select()
